$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 59, shifting rows 59:86 down to 60:87
$ws.Rows.Item(59).Insert()

# Populate the new row 59 with fresh data (values copied from row 60 for the
# fields that are unchanged, per the diff)
$ws.Cells.Item(59, 1).Value = 2
$ws.Cells.Item(59, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = (Get-Date -Year 2021 -Month 10 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100112024
$ws.Cells.Item(59, 7).Value = "Choclo"
$ws.Cells.Item(59, 8).Value = "Dulce o Americano"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 500
$ws.Cells.Item(59, 11).Value = 39000
$ws.Cells.Item(59, 12).Value = 40000
$ws.Cells.Item(59, 13).Value = 39500
$ws.Cells.Item(59, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(59, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(59, 16).Value = 564
$ws.Cells.Item(59, 17).Value = 70
$ws.Cells.Item(59, 18).Value = "Hortaliza"
